$wb = $excel.ActiveWorkbook

# Duplicate the last sheet ("Лист1") - the copy is placed immediately
# before it, becomes the active sheet, and inherits the next sheetId.
# This reproduces Excel's native "Move or Copy... (Create a copy)" behaviour,
# matching the workbook.xml <sheets>/<sheetId>/<r:id> renumbering in the diff.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($lastSheet)
$newSheet = $wb.ActiveSheet

# Paint the new sample (digit "2") pixel grid, A1:T20, 0/1 values.
$newSheet.Cells.Item(1,1).Value = 0
$newSheet.Cells.Item(1,2).Value = 0
$newSheet.Cells.Item(1,3).Value = 0
$newSheet.Cells.Item(1,4).Value = 0
$newSheet.Cells.Item(1,5).Value = 0
$newSheet.Cells.Item(1,6).Value = 0
$newSheet.Cells.Item(1,7).Value = 0
$newSheet.Cells.Item(1,8).Value = 0
$newSheet.Cells.Item(1,9).Value = 0
$newSheet.Cells.Item(1,10).Value = 0
$newSheet.Cells.Item(1,11).Value = 0
$newSheet.Cells.Item(1,12).Value = 0
$newSheet.Cells.Item(1,13).Value = 0
$newSheet.Cells.Item(1,14).Value = 0
$newSheet.Cells.Item(1,15).Value = 0
$newSheet.Cells.Item(1,16).Value = 0
$newSheet.Cells.Item(1,17).Value = 0
$newSheet.Cells.Item(1,18).Value = 0
$newSheet.Cells.Item(1,19).Value = 0
$newSheet.Cells.Item(1,20).Value = 0
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = 0
$newSheet.Cells.Item(2,3).Value = 0
$newSheet.Cells.Item(2,4).Value = 0
$newSheet.Cells.Item(2,5).Value = 0
$newSheet.Cells.Item(2,6).Value = 0
$newSheet.Cells.Item(2,7).Value = 0
$newSheet.Cells.Item(2,8).Value = 0
$newSheet.Cells.Item(2,9).Value = 0
$newSheet.Cells.Item(2,10).Value = 0
$newSheet.Cells.Item(2,11).Value = 0
$newSheet.Cells.Item(2,12).Value = 0
$newSheet.Cells.Item(2,13).Value = 0
$newSheet.Cells.Item(2,14).Value = 0
$newSheet.Cells.Item(2,15).Value = 0
$newSheet.Cells.Item(2,16).Value = 0
$newSheet.Cells.Item(2,17).Value = 0
$newSheet.Cells.Item(2,18).Value = 0
$newSheet.Cells.Item(2,19).Value = 0
$newSheet.Cells.Item(2,20).Value = 0
$newSheet.Cells.Item(3,1).Value = 0
$newSheet.Cells.Item(3,2).Value = 0
$newSheet.Cells.Item(3,3).Value = 0
$newSheet.Cells.Item(3,4).Value = 0
$newSheet.Cells.Item(3,5).Value = 0
$newSheet.Cells.Item(3,6).Value = 0
$newSheet.Cells.Item(3,7).Value = 0
$newSheet.Cells.Item(3,8).Value = 0
$newSheet.Cells.Item(3,9).Value = 0
$newSheet.Cells.Item(3,10).Value = 0
$newSheet.Cells.Item(3,11).Value = 0
$newSheet.Cells.Item(3,12).Value = 0
$newSheet.Cells.Item(3,13).Value = 0
$newSheet.Cells.Item(3,14).Value = 0
$newSheet.Cells.Item(3,15).Value = 0
$newSheet.Cells.Item(3,16).Value = 0
$newSheet.Cells.Item(3,17).Value = 0
$newSheet.Cells.Item(3,18).Value = 0
$newSheet.Cells.Item(3,19).Value = 0
$newSheet.Cells.Item(3,20).Value = 0
$newSheet.Cells.Item(4,1).Value = 0
$newSheet.Cells.Item(4,2).Value = 0
$newSheet.Cells.Item(4,3).Value = 0
$newSheet.Cells.Item(4,4).Value = 0
$newSheet.Cells.Item(4,5).Value = 0
$newSheet.Cells.Item(4,6).Value = 0
$newSheet.Cells.Item(4,7).Value = 0
$newSheet.Cells.Item(4,8).Value = 1
$newSheet.Cells.Item(4,9).Value = 1
$newSheet.Cells.Item(4,10).Value = 1
$newSheet.Cells.Item(4,11).Value = 1
$newSheet.Cells.Item(4,12).Value = 1
$newSheet.Cells.Item(4,13).Value = 1
$newSheet.Cells.Item(4,14).Value = 0
$newSheet.Cells.Item(4,15).Value = 0
$newSheet.Cells.Item(4,16).Value = 0
$newSheet.Cells.Item(4,17).Value = 0
$newSheet.Cells.Item(4,18).Value = 0
$newSheet.Cells.Item(4,19).Value = 0
$newSheet.Cells.Item(4,20).Value = 0
$newSheet.Cells.Item(5,1).Value = 0
$newSheet.Cells.Item(5,2).Value = 0
$newSheet.Cells.Item(5,3).Value = 0
$newSheet.Cells.Item(5,4).Value = 0
$newSheet.Cells.Item(5,5).Value = 0
$newSheet.Cells.Item(5,6).Value = 0
$newSheet.Cells.Item(5,7).Value = 0
$newSheet.Cells.Item(5,8).Value = 1
$newSheet.Cells.Item(5,9).Value = 0
$newSheet.Cells.Item(5,10).Value = 0
$newSheet.Cells.Item(5,11).Value = 0
$newSheet.Cells.Item(5,12).Value = 0
$newSheet.Cells.Item(5,13).Value = 1
$newSheet.Cells.Item(5,14).Value = 1
$newSheet.Cells.Item(5,15).Value = 0
$newSheet.Cells.Item(5,16).Value = 0
$newSheet.Cells.Item(5,17).Value = 0
$newSheet.Cells.Item(5,18).Value = 0
$newSheet.Cells.Item(5,19).Value = 0
$newSheet.Cells.Item(5,20).Value = 0
$newSheet.Cells.Item(6,1).Value = 0
$newSheet.Cells.Item(6,2).Value = 0
$newSheet.Cells.Item(6,3).Value = 0
$newSheet.Cells.Item(6,4).Value = 0
$newSheet.Cells.Item(6,5).Value = 0
$newSheet.Cells.Item(6,6).Value = 0
$newSheet.Cells.Item(6,7).Value = 0
$newSheet.Cells.Item(6,8).Value = 0
$newSheet.Cells.Item(6,9).Value = 0
$newSheet.Cells.Item(6,10).Value = 0
$newSheet.Cells.Item(6,11).Value = 0
$newSheet.Cells.Item(6,12).Value = 0
$newSheet.Cells.Item(6,13).Value = 0
$newSheet.Cells.Item(6,14).Value = 1
$newSheet.Cells.Item(6,15).Value = 0
$newSheet.Cells.Item(6,16).Value = 0
$newSheet.Cells.Item(6,17).Value = 0
$newSheet.Cells.Item(6,18).Value = 0
$newSheet.Cells.Item(6,19).Value = 0
$newSheet.Cells.Item(6,20).Value = 0
$newSheet.Cells.Item(7,1).Value = 0
$newSheet.Cells.Item(7,2).Value = 0
$newSheet.Cells.Item(7,3).Value = 0
$newSheet.Cells.Item(7,4).Value = 0
$newSheet.Cells.Item(7,5).Value = 0
$newSheet.Cells.Item(7,6).Value = 0
$newSheet.Cells.Item(7,7).Value = 0
$newSheet.Cells.Item(7,8).Value = 0
$newSheet.Cells.Item(7,9).Value = 0
$newSheet.Cells.Item(7,10).Value = 0
$newSheet.Cells.Item(7,11).Value = 0
$newSheet.Cells.Item(7,12).Value = 0
$newSheet.Cells.Item(7,13).Value = 0
$newSheet.Cells.Item(7,14).Value = 1
$newSheet.Cells.Item(7,15).Value = 0
$newSheet.Cells.Item(7,16).Value = 0
$newSheet.Cells.Item(7,17).Value = 0
$newSheet.Cells.Item(7,18).Value = 0
$newSheet.Cells.Item(7,19).Value = 0
$newSheet.Cells.Item(7,20).Value = 0
$newSheet.Cells.Item(8,1).Value = 0
$newSheet.Cells.Item(8,2).Value = 0
$newSheet.Cells.Item(8,3).Value = 0
$newSheet.Cells.Item(8,4).Value = 0
$newSheet.Cells.Item(8,5).Value = 0
$newSheet.Cells.Item(8,6).Value = 0
$newSheet.Cells.Item(8,7).Value = 0
$newSheet.Cells.Item(8,8).Value = 0
$newSheet.Cells.Item(8,9).Value = 0
$newSheet.Cells.Item(8,10).Value = 0
$newSheet.Cells.Item(8,11).Value = 0
$newSheet.Cells.Item(8,12).Value = 0
$newSheet.Cells.Item(8,13).Value = 0
$newSheet.Cells.Item(8,14).Value = 1
$newSheet.Cells.Item(8,15).Value = 0
$newSheet.Cells.Item(8,16).Value = 0
$newSheet.Cells.Item(8,17).Value = 0
$newSheet.Cells.Item(8,18).Value = 0
$newSheet.Cells.Item(8,19).Value = 0
$newSheet.Cells.Item(8,20).Value = 0
$newSheet.Cells.Item(9,1).Value = 0
$newSheet.Cells.Item(9,2).Value = 0
$newSheet.Cells.Item(9,3).Value = 0
$newSheet.Cells.Item(9,4).Value = 0
$newSheet.Cells.Item(9,5).Value = 0
$newSheet.Cells.Item(9,6).Value = 0
$newSheet.Cells.Item(9,7).Value = 0
$newSheet.Cells.Item(9,8).Value = 0
$newSheet.Cells.Item(9,9).Value = 0
$newSheet.Cells.Item(9,10).Value = 0
$newSheet.Cells.Item(9,11).Value = 0
$newSheet.Cells.Item(9,12).Value = 0
$newSheet.Cells.Item(9,13).Value = 0
$newSheet.Cells.Item(9,14).Value = 1
$newSheet.Cells.Item(9,15).Value = 0
$newSheet.Cells.Item(9,16).Value = 0
$newSheet.Cells.Item(9,17).Value = 0
$newSheet.Cells.Item(9,18).Value = 0
$newSheet.Cells.Item(9,19).Value = 0
$newSheet.Cells.Item(9,20).Value = 0
$newSheet.Cells.Item(10,1).Value = 0
$newSheet.Cells.Item(10,2).Value = 0
$newSheet.Cells.Item(10,3).Value = 0
$newSheet.Cells.Item(10,4).Value = 0
$newSheet.Cells.Item(10,5).Value = 0
$newSheet.Cells.Item(10,6).Value = 0
$newSheet.Cells.Item(10,7).Value = 0
$newSheet.Cells.Item(10,8).Value = 0
$newSheet.Cells.Item(10,9).Value = 0
$newSheet.Cells.Item(10,10).Value = 0
$newSheet.Cells.Item(10,11).Value = 0
$newSheet.Cells.Item(10,12).Value = 0
$newSheet.Cells.Item(10,13).Value = 0
$newSheet.Cells.Item(10,14).Value = 1
$newSheet.Cells.Item(10,15).Value = 0
$newSheet.Cells.Item(10,16).Value = 0
$newSheet.Cells.Item(10,17).Value = 0
$newSheet.Cells.Item(10,18).Value = 0
$newSheet.Cells.Item(10,19).Value = 0
$newSheet.Cells.Item(10,20).Value = 0
$newSheet.Cells.Item(11,1).Value = 0
$newSheet.Cells.Item(11,2).Value = 0
$newSheet.Cells.Item(11,3).Value = 0
$newSheet.Cells.Item(11,4).Value = 0
$newSheet.Cells.Item(11,5).Value = 0
$newSheet.Cells.Item(11,6).Value = 0
$newSheet.Cells.Item(11,7).Value = 0
$newSheet.Cells.Item(11,8).Value = 0
$newSheet.Cells.Item(11,9).Value = 0
$newSheet.Cells.Item(11,10).Value = 0
$newSheet.Cells.Item(11,11).Value = 0
$newSheet.Cells.Item(11,12).Value = 0
$newSheet.Cells.Item(11,13).Value = 1
$newSheet.Cells.Item(11,14).Value = 1
$newSheet.Cells.Item(11,15).Value = 0
$newSheet.Cells.Item(11,16).Value = 0
$newSheet.Cells.Item(11,17).Value = 0
$newSheet.Cells.Item(11,18).Value = 0
$newSheet.Cells.Item(11,19).Value = 0
$newSheet.Cells.Item(11,20).Value = 0
$newSheet.Cells.Item(12,1).Value = 0
$newSheet.Cells.Item(12,2).Value = 0
$newSheet.Cells.Item(12,3).Value = 0
$newSheet.Cells.Item(12,4).Value = 0
$newSheet.Cells.Item(12,5).Value = 0
$newSheet.Cells.Item(12,6).Value = 0
$newSheet.Cells.Item(12,7).Value = 0
$newSheet.Cells.Item(12,8).Value = 0
$newSheet.Cells.Item(12,9).Value = 0
$newSheet.Cells.Item(12,10).Value = 0
$newSheet.Cells.Item(12,11).Value = 0
$newSheet.Cells.Item(12,12).Value = 1
$newSheet.Cells.Item(12,13).Value = 1
$newSheet.Cells.Item(12,14).Value = 0
$newSheet.Cells.Item(12,15).Value = 0
$newSheet.Cells.Item(12,16).Value = 0
$newSheet.Cells.Item(12,17).Value = 0
$newSheet.Cells.Item(12,18).Value = 0
$newSheet.Cells.Item(12,19).Value = 0
$newSheet.Cells.Item(12,20).Value = 0
$newSheet.Cells.Item(13,1).Value = 0
$newSheet.Cells.Item(13,2).Value = 0
$newSheet.Cells.Item(13,3).Value = 0
$newSheet.Cells.Item(13,4).Value = 0
$newSheet.Cells.Item(13,5).Value = 0
$newSheet.Cells.Item(13,6).Value = 0
$newSheet.Cells.Item(13,7).Value = 0
$newSheet.Cells.Item(13,8).Value = 0
$newSheet.Cells.Item(13,9).Value = 0
$newSheet.Cells.Item(13,10).Value = 0
$newSheet.Cells.Item(13,11).Value = 1
$newSheet.Cells.Item(13,12).Value = 1
$newSheet.Cells.Item(13,13).Value = 0
$newSheet.Cells.Item(13,14).Value = 0
$newSheet.Cells.Item(13,15).Value = 0
$newSheet.Cells.Item(13,16).Value = 0
$newSheet.Cells.Item(13,17).Value = 0
$newSheet.Cells.Item(13,18).Value = 0
$newSheet.Cells.Item(13,19).Value = 0
$newSheet.Cells.Item(13,20).Value = 0
$newSheet.Cells.Item(14,1).Value = 0
$newSheet.Cells.Item(14,2).Value = 0
$newSheet.Cells.Item(14,3).Value = 0
$newSheet.Cells.Item(14,4).Value = 0
$newSheet.Cells.Item(14,5).Value = 0
$newSheet.Cells.Item(14,6).Value = 0
$newSheet.Cells.Item(14,7).Value = 0
$newSheet.Cells.Item(14,8).Value = 0
$newSheet.Cells.Item(14,9).Value = 0
$newSheet.Cells.Item(14,10).Value = 1
$newSheet.Cells.Item(14,11).Value = 1
$newSheet.Cells.Item(14,12).Value = 0
$newSheet.Cells.Item(14,13).Value = 0
$newSheet.Cells.Item(14,14).Value = 0
$newSheet.Cells.Item(14,15).Value = 0
$newSheet.Cells.Item(14,16).Value = 0
$newSheet.Cells.Item(14,17).Value = 0
$newSheet.Cells.Item(14,18).Value = 0
$newSheet.Cells.Item(14,19).Value = 0
$newSheet.Cells.Item(14,20).Value = 0
$newSheet.Cells.Item(15,1).Value = 0
$newSheet.Cells.Item(15,2).Value = 0
$newSheet.Cells.Item(15,3).Value = 0
$newSheet.Cells.Item(15,4).Value = 0
$newSheet.Cells.Item(15,5).Value = 0
$newSheet.Cells.Item(15,6).Value = 0
$newSheet.Cells.Item(15,7).Value = 0
$newSheet.Cells.Item(15,8).Value = 1
$newSheet.Cells.Item(15,9).Value = 1
$newSheet.Cells.Item(15,10).Value = 1
$newSheet.Cells.Item(15,11).Value = 0
$newSheet.Cells.Item(15,12).Value = 0
$newSheet.Cells.Item(15,13).Value = 0
$newSheet.Cells.Item(15,14).Value = 0
$newSheet.Cells.Item(15,15).Value = 0
$newSheet.Cells.Item(15,16).Value = 0
$newSheet.Cells.Item(15,17).Value = 0
$newSheet.Cells.Item(15,18).Value = 0
$newSheet.Cells.Item(15,19).Value = 0
$newSheet.Cells.Item(15,20).Value = 0
$newSheet.Cells.Item(16,1).Value = 0
$newSheet.Cells.Item(16,2).Value = 0
$newSheet.Cells.Item(16,3).Value = 0
$newSheet.Cells.Item(16,4).Value = 0
$newSheet.Cells.Item(16,5).Value = 0
$newSheet.Cells.Item(16,6).Value = 1
$newSheet.Cells.Item(16,7).Value = 1
$newSheet.Cells.Item(16,8).Value = 1
$newSheet.Cells.Item(16,9).Value = 1
$newSheet.Cells.Item(16,10).Value = 1
$newSheet.Cells.Item(16,11).Value = 0
$newSheet.Cells.Item(16,12).Value = 0
$newSheet.Cells.Item(16,13).Value = 0
$newSheet.Cells.Item(16,14).Value = 0
$newSheet.Cells.Item(16,15).Value = 0
$newSheet.Cells.Item(16,16).Value = 0
$newSheet.Cells.Item(16,17).Value = 0
$newSheet.Cells.Item(16,18).Value = 0
$newSheet.Cells.Item(16,19).Value = 0
$newSheet.Cells.Item(16,20).Value = 0
$newSheet.Cells.Item(17,1).Value = 0
$newSheet.Cells.Item(17,2).Value = 0
$newSheet.Cells.Item(17,3).Value = 0
$newSheet.Cells.Item(17,4).Value = 0
$newSheet.Cells.Item(17,5).Value = 1
$newSheet.Cells.Item(17,6).Value = 1
$newSheet.Cells.Item(17,7).Value = 0
$newSheet.Cells.Item(17,8).Value = 0
$newSheet.Cells.Item(17,9).Value = 0
$newSheet.Cells.Item(17,10).Value = 1
$newSheet.Cells.Item(17,11).Value = 1
$newSheet.Cells.Item(17,12).Value = 0
$newSheet.Cells.Item(17,13).Value = 0
$newSheet.Cells.Item(17,14).Value = 0
$newSheet.Cells.Item(17,15).Value = 0
$newSheet.Cells.Item(17,16).Value = 1
$newSheet.Cells.Item(17,17).Value = 0
$newSheet.Cells.Item(17,18).Value = 0
$newSheet.Cells.Item(17,19).Value = 0
$newSheet.Cells.Item(17,20).Value = 0
$newSheet.Cells.Item(18,1).Value = 0
$newSheet.Cells.Item(18,2).Value = 0
$newSheet.Cells.Item(18,3).Value = 0
$newSheet.Cells.Item(18,4).Value = 0
$newSheet.Cells.Item(18,5).Value = 0
$newSheet.Cells.Item(18,6).Value = 0
$newSheet.Cells.Item(18,7).Value = 0
$newSheet.Cells.Item(18,8).Value = 0
$newSheet.Cells.Item(18,9).Value = 0
$newSheet.Cells.Item(18,10).Value = 0
$newSheet.Cells.Item(18,11).Value = 1
$newSheet.Cells.Item(18,12).Value = 1
$newSheet.Cells.Item(18,13).Value = 1
$newSheet.Cells.Item(18,14).Value = 1
$newSheet.Cells.Item(18,15).Value = 1
$newSheet.Cells.Item(18,16).Value = 1
$newSheet.Cells.Item(18,17).Value = 0
$newSheet.Cells.Item(18,18).Value = 0
$newSheet.Cells.Item(18,19).Value = 0
$newSheet.Cells.Item(18,20).Value = 0
$newSheet.Cells.Item(19,1).Value = 0
$newSheet.Cells.Item(19,2).Value = 0
$newSheet.Cells.Item(19,3).Value = 0
$newSheet.Cells.Item(19,4).Value = 0
$newSheet.Cells.Item(19,5).Value = 0
$newSheet.Cells.Item(19,6).Value = 0
$newSheet.Cells.Item(19,7).Value = 0
$newSheet.Cells.Item(19,8).Value = 0
$newSheet.Cells.Item(19,9).Value = 0
$newSheet.Cells.Item(19,10).Value = 0
$newSheet.Cells.Item(19,11).Value = 0
$newSheet.Cells.Item(19,12).Value = 0
$newSheet.Cells.Item(19,13).Value = 0
$newSheet.Cells.Item(19,14).Value = 0
$newSheet.Cells.Item(19,15).Value = 0
$newSheet.Cells.Item(19,16).Value = 0
$newSheet.Cells.Item(19,17).Value = 0
$newSheet.Cells.Item(19,18).Value = 0
$newSheet.Cells.Item(19,19).Value = 0
$newSheet.Cells.Item(19,20).Value = 0
$newSheet.Cells.Item(20,1).Value = 0
$newSheet.Cells.Item(20,2).Value = 0
$newSheet.Cells.Item(20,3).Value = 0
$newSheet.Cells.Item(20,4).Value = 0
$newSheet.Cells.Item(20,5).Value = 0
$newSheet.Cells.Item(20,6).Value = 0
$newSheet.Cells.Item(20,7).Value = 0
$newSheet.Cells.Item(20,8).Value = 0
$newSheet.Cells.Item(20,9).Value = 0
$newSheet.Cells.Item(20,10).Value = 0
$newSheet.Cells.Item(20,11).Value = 0
$newSheet.Cells.Item(20,12).Value = 0
$newSheet.Cells.Item(20,13).Value = 0
$newSheet.Cells.Item(20,14).Value = 0
$newSheet.Cells.Item(20,15).Value = 0
$newSheet.Cells.Item(20,16).Value = 0
$newSheet.Cells.Item(20,17).Value = 0
$newSheet.Cells.Item(20,18).Value = 0
$newSheet.Cells.Item(20,19).Value = 0
$newSheet.Cells.Item(20,20).Value = 0

# Keep the new sheet as the active / selected tab (matches activeTab="9").
$newSheet.Select()
